# Insert one new weekly price record at row 434, pushing existing data
# (old rows 434-531) down by one row (new rows 435-532).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 434; Excel shifts rows 434..531
# down to 435..532 and the new (currently blank) row becomes row 434.
$ws.Rows.Item(434).Insert()

# Populate the newly inserted row 434 with the new weekly record.
$ws.Cells.Item(434, 1).Value2 = 5
$ws.Cells.Item(434, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(434, 3).Value2 = "Maule"
$ws.Cells.Item(434, 4).Value2 = 45244
$ws.Cells.Item(434, 5).Value2 = 7
$ws.Cells.Item(434, 6).Value2 = 100112009
$ws.Cells.Item(434, 7).Value2 = "Acelga"
$ws.Cells.Item(434, 8).Value2 = "Sin especificar"
$ws.Cells.Item(434, 9).Value2 = "Primera"
$ws.Cells.Item(434, 10).Value2 = 500
$ws.Cells.Item(434, 11).Value2 = 2300
$ws.Cells.Item(434, 12).Value2 = 2300
$ws.Cells.Item(434, 13).Value2 = 2300
$ws.Cells.Item(434, 14).Value2 = "$/docena de atados (4 kilos)"
$ws.Cells.Item(434, 15).Value2 = "Región del Maule"
$ws.Cells.Item(434, 16).Value2 = 575
$ws.Cells.Item(434, 17).Value2 = 4
$ws.Cells.Item(434, 18).Value2 = "Hortaliza"
